# Översikt GRUMS.xlsx - automatic update of files
# 1) Column C ("Förändrad") date bumped from 45184 to 45186 for every data row.
# 2) For the first four data rows (2-5), the HYPERLINK() formulas in columns
#    S, T, V, W, X, Y gain a second "friendly name" argument equal to the
#    row's "Beteckning" (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 270
$oldDate = 45184
$newDate = 45186

# --- 1) Bump the "Förändrad" date in column C for every data row ---------
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value = $newDate
    }
}

# --- 2) Add friendly display text to the HYPERLINK formulas in rows 2-5 --
$linkColumns = @(
    @{ Col = 19; Folder = "artfynd";         Ext = "xlsx" },  # S
    @{ Col = 20; Folder = "kartor";          Ext = "png"  },  # T
    @{ Col = 22; Folder = "klagomål";        Ext = "docx" },  # V
    @{ Col = 23; Folder = "klagomålsmail";   Ext = "docx" },  # W
    @{ Col = 24; Folder = "tillsyn";         Ext = "docx" },  # X
    @{ Col = 25; Folder = "tillsynsmail";    Ext = "docx" }   # Y
)

for ($r = 2; $r -le 5; $r++) {
    $designation = $ws.Cells.Item($r, 1).Text

    foreach ($link in $linkColumns) {
        $url = "https://klasma.github.io/Logging_GRUMS/" + $link.Folder + "/" + $designation + "." + $link.Ext
        $formula = '=HYPERLINK("' + $url + '", "' + $designation + '")'
        $ws.Cells.Item($r, $link.Col).Formula = $formula
    }
}
